$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.821486949920654
$ws.Range("B1").Value = 4.535807609558105
$ws.Range("C1").Value = 4.132960319519043
$ws.Range("D1").Value = 0.9058610200881958
$ws.Range("E1").Value = 0.4753615856170654
